$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for D,E,F,G,H (open_price, close_price, high_price, low_price, shares_outstanding)
# and I (fixed_ticker) for rows 2-7. All fixed_ticker values now point to "CRWV" (same as ticker in column A).

$data = @(
    @{ Row = 2; D = 39;    E = 37.08000183105469; F = 41.93999862670898; G = 36;            H = 380162985 },
    @{ Row = 3; D = 39;    E = 37.08000183105469; F = 41.93999862670898; G = 36;            H = 380162985 },
    @{ Row = 4; D = 39;    E = 37.08000183105469; F = 41.93999862670898; G = 36;            H = 380162985 },
    @{ Row = 5; D = 39;    E = 37.08000183105469; F = 41.93999862670898; G = 36;            H = 380162985 },
    @{ Row = 6; D = 38.5;  E = 41.29999923706055; F = 64.62000274658203; G = 33.51499938964844; H = 380162985 },
    @{ Row = 7; D = 162.7799987792969; E = 114.129997253418; F = 166.2200012207031; G = 100.8000030517578; H = 380162985 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
    $ws.Range("H$r").Value = $entry.H
    $ws.Range("I$r").Value = "CRWV"
}
